# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2404"
#   "<header>_new" -> "<header>_FV2410"
# Then turn the header row + data range into a native Excel Table, and
# freeze the header row (row 1) in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row cells (A1:U1) ------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2404"
    "B1" = "Segmentgruppe_FV2404"
    "C1" = "Segment_FV2404"
    "D1" = "Datenelement_FV2404"
    "E1" = "Segment ID_FV2404"
    "F1" = "Code_FV2404"
    "G1" = "Qualifier_FV2404"
    "H1" = "Beschreibung_FV2404"
    "I1" = "Bedingungsausdruck_FV2404"
    "J1" = "Bedingung_FV2404"
    "K1" = "diff"
    "L1" = "Segmentname_FV2410"
    "M1" = "Segmentgruppe_FV2410"
    "N1" = "Segment_FV2410"
    "O1" = "Datenelement_FV2410"
    "P1" = "Segment ID_FV2410"
    "Q1" = "Code_FV2410"
    "R1" = "Qualifier_FV2410"
    "S1" = "Beschreibung_FV2410"
    "T1" = "Bedingungsausdruck_FV2410"
    "U1" = "Bedingung_FV2410"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2) Turn A1:U80 into a real Excel Table (ListObject) -----------------
$dataRange = $ws.Range("A1:U80")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"

# --- 3) Freeze the header row (row 1) ------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
